$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CoinCell {
    param($cell, $value)
    if ($value -match '^-?\d+(\.\d+)?$') {
        $ws.Range($cell).Value = "'" + $value
    } else {
        $ws.Range($cell).Value = $value
    }
}

# Row 2
Set-CoinCell "D2" '65.979.03'
Set-CoinCell "E2" '  -1.69%  '

# Row 3
Set-CoinCell "D3" '3.484.36'
Set-CoinCell "E3" '  -2.44%  '

# Row 4
Set-CoinCell "E4" '  -0.06%  '

# Row 5
Set-CoinCell "D5" '582.69'
Set-CoinCell "E5" '  +5.45%  '

# Row 6
Set-CoinCell "D6" '177.79'
Set-CoinCell "E6" '  -5.95%  '

# Row 7
Set-CoinCell "D7" '0.632'
Set-CoinCell "E7" '  +3.74%  '

# Row 8
Set-CoinCell "E8" '  -0.07%  '

# Row 9
Set-CoinCell "E9" '  -0.20%  '

# Row 10
Set-CoinCell "E10" '  +5.01%  '

# Row 11
Set-CoinCell "D11" '55.89'
Set-CoinCell "E11" '  +1.75%  '

# Row 12
Set-CoinCell "E12" '  +2.16%  '

# Row 13
Set-CoinCell "D13" '9.24'
Set-CoinCell "E13" '  -2.10%  '

# Row 14
Set-CoinCell "D14" '4.041.54'
Set-CoinCell "E14" '  -2.32%  '

# Row 15
Set-CoinCell "D15" '3.483.11'
Set-CoinCell "E15" '  -2.57%  '

# Row 16
Set-CoinCell "E16" '  -0.06%  '

# Row 17
Set-CoinCell "D17" '18.26'
Set-CoinCell "E17" '  -0.11%  '

# Row 18
Set-CoinCell "B18" 'Uniswap'
Set-CoinCell "C18" 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-CoinCell "D18" '12.06'
Set-CoinCell "E18" '  +0.60%  '

# Row 19
Set-CoinCell "B19" 'WrappedBTC'
Set-CoinCell "C19" 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-CoinCell "D19" '65.926.24'
Set-CoinCell "E19" '  -1.82%  '

# Row 20
Set-CoinCell "E20" '  +1.37%  '

# Row 21
Set-CoinCell "D21" '411.30'
Set-CoinCell "E21" '  -5.10%  '

# Row 22
Set-CoinCell "D22" '4.29'
Set-CoinCell "E22" '  +9.43%  '

# Row 23
Set-CoinCell "D23" '4.38'
Set-CoinCell "E23" '  +5.87%  '

# Row 24
Set-CoinCell "D24" '84.70'
Set-CoinCell "E24" '  -1.09%  '

# Row 25
Set-CoinCell "D25" '13.42'
Set-CoinCell "E25" '  +10.65%  '

# Row 26
Set-CoinCell "D26" '11.05'
Set-CoinCell "E26" '  -0.54%  '

# Row 27
Set-CoinCell "D27" '2.86'
Set-CoinCell "E27" '  -1.53%  '

# Row 28
Set-CoinCell "D28" '6.04'
Set-CoinCell "E28" '  -0.21%  '

# Row 29
Set-CoinCell "D29" '9.18'
Set-CoinCell "E29" '  +1.71%  '

# Row 30
Set-CoinCell "D30" '30.20'
Set-CoinCell "E30" '  -1.05%  '

# Row 31
Set-CoinCell "E31" '  +0.34%  '

# Row 32
Set-CoinCell "D32" '11.73'
Set-CoinCell "E32" '  -0.33%  '

# Row 33
Set-CoinCell "D33" '593.13'
Set-CoinCell "E33" '  -8.64%  '

# Row 34
Set-CoinCell "E34" '  -1.49%  '

# Row 35
Set-CoinCell "D35" '60.86'
Set-CoinCell "E35" '  +2.00%  '

# Row 36
Set-CoinCell "D36" '0.153'
Set-CoinCell "E36" '  +1.12%  '

# Row 37
Set-CoinCell "E37" '  +0.08%  '

# Row 38
Set-CoinCell "B38" 'PEPE'
Set-CoinCell "C38" 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-CoinCell "D38" '0.0₃0794'
Set-CoinCell "E38" '  -3.75%  '

# Row 39
Set-CoinCell "B39" 'Stacks'
Set-CoinCell "C39" 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-CoinCell "D39" '3.59'
Set-CoinCell "E39" '  +7.10%  '

# Row 40
Set-CoinCell "D40" '36.83'
Set-CoinCell "E40" '  -4.77%  '

# Row 41
Set-CoinCell "D41" '0.384'
Set-CoinCell "E41" '  -2.03%  '

# Row 42
Set-CoinCell "D42" '3.222.47'
Set-CoinCell "E42" '  +5.69%  '

# Row 43
Set-CoinCell "D43" '0.999'
Set-CoinCell "E43" '  -0.08%  '

# Row 44
Set-CoinCell "D44" '2.96'
Set-CoinCell "E44" '  +2.27%  '

# Row 45
Set-CoinCell "B45" 'Fetch.AI'
Set-CoinCell "C45" 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-CoinCell "D45" '2.55'
Set-CoinCell "E45" '  -4.91%  '

# Row 46
Set-CoinCell "B46" 'ApeXProtocol'
Set-CoinCell "C46" 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-CoinCell "D46" '3.31'
Set-CoinCell "E46" '  -3.05%  '

# Row 47
Set-CoinCell "E47" '  -0.52%  '

# Row 48
Set-CoinCell "B48" 'Stellar'
Set-CoinCell "C48" 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-CoinCell "D48" '0.133'
Set-CoinCell "E48" '  +1.34%  '

# Row 49
Set-CoinCell "B49" 'WEMIXToken'
Set-CoinCell "C49" 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-CoinCell "D49" '2.65'
Set-CoinCell "E49" '  -5.00%  '

# Row 50
Set-CoinCell "D50" '8.59'
Set-CoinCell "E50" '  -1.64%  '

# Row 51
Set-CoinCell "D51" '139.69'
Set-CoinCell "E51" '  -1.47%  '
